$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR EXPENSES")

# Department labels used in both the TOTAL table (rows 5-11)
# and the AVERAGE table (rows 15-21)
$departments = @(
    "Accounting",
    "Customer Service",
    "Executive",
    "Human Resources",
    "Facilities",
    "IT",
    "Sales"
)

for ($i = 0; $i -lt $departments.Length; $i++) {
    $row1 = 5 + $i
    $row2 = 15 + $i
    $ws.Range("B$row1").Value = $departments[$i]
    $ws.Range("B$row2").Value = $departments[$i]
}

# Quarterly leave-days totals (first "HR LEAVE DAYS" table, rows 5-11)
$totals = @(
    @(31, 28, 21, 19),
    @(21, 15, 25, 41),
    @(10, 6, 12, 4),
    @(20, 31, 38, 29),
    @(1, 0, 5, 0),
    @(25, 33, 22, 37),
    @(36, 11, 40, 56)
)

for ($i = 0; $i -lt $totals.Length; $i++) {
    $row = 5 + $i
    $ws.Range("C$row").Value = $totals[$i][0]
    $ws.Range("D$row").Value = $totals[$i][1]
    $ws.Range("E$row").Value = $totals[$i][2]
    $ws.Range("F$row").Value = $totals[$i][3]
}

# Average leave days per month (consolidated "SUMMARY" table, rows 15-21)
$averages = @(
    @(10.333333333333334, 9.3333333333333339, 7, 6.333333333333333),
    @(7, 5, 8.3333333333333339, 13.666666666666666),
    @(5, 3, 6, 2),
    @(6.666666666666667, 10.333333333333334, 12.666666666666666, 9.6666666666666661),
    @(1, 0, 5, 0),
    @(8.3333333333333339, 11, 7.333333333333333, 12.333333333333334),
    @(12, 3.6666666666666665, 13.333333333333334, 18.666666666666668)
)

for ($i = 0; $i -lt $averages.Length; $i++) {
    $row = 15 + $i
    $ws.Range("C$row").Value = $averages[$i][0]
    $ws.Range("D$row").Value = $averages[$i][1]
    $ws.Range("E$row").Value = $averages[$i][2]
    $ws.Range("F$row").Value = $averages[$i][3]
}

# Reflect the view scrolled down to the SUMMARY table with B15:F21 selected
$ws.Range("B15:F21").Select()
$excel.ActiveWindow.ScrollRow = 13
